$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (PRU)
$ws.Range("D2").Value = 108.25
$ws.Range("F2").Value = 3.97
$ws.Range("N2").Value = 85.83574689470727

# Row 3 (UNH)
$ws.Range("D3").Value = 329.77
$ws.Range("F3").Value = 5.85
$ws.Range("N3").Value = 85.83574689470727

# Row 4 (MET) - only MACRO_SCORE changes
$ws.Range("N4").Value = 85.83574689470727

# Row 5 (AIG)
$ws.Range("D5").Value = 76.16
$ws.Range("F5").Value = 0.62
$ws.Range("N5").Value = 85.83574689470727
